$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "V-1769915744936"
$ws.Range("C2").Value = "10:15 p. m."
$ws.Range("E2").Value = "Aguardiente Amarillo Botella (x1)"
$ws.Range("F2").Value = 108000

# Add row 3
$ws.Range("A3").Value = "V-1769915748653"
$ws.Range("B3").Value = "31/1/2026"
$ws.Range("C3").Value = "10:15 p. m."
$ws.Range("D3").Value = "Martha"
$ws.Range("E3").Value = "Aguardiente Amarillo Botella (x1)"
$ws.Range("F3").Value = 108000
$ws.Range("G3").Value = 0

# Add row 4
$ws.Range("A4").Value = "V-1769915756171"
$ws.Range("B4").Value = "31/1/2026"
$ws.Range("C4").Value = "10:15 p. m."
$ws.Range("D4").Value = "Martha"
$ws.Range("E4").Value = "Cerveza Corona (x4)"
$ws.Range("F4").Value = 40000
$ws.Range("G4").Value = 0

# Add row 5
$ws.Range("A5").Value = "V-1769915763343"
$ws.Range("B5").Value = "31/1/2026"
$ws.Range("C5").Value = "10:16 p. m."
$ws.Range("D5").Value = "Martha"
$ws.Range("E5").Value = "Aguardiente Amarillo Botella (x1)"
$ws.Range("F5").Value = 108000
$ws.Range("G5").Value = 1
